$d = $word.ActiveDocument

# --- Change 1: merge the "outline in RMarkdown..." paragraph runs (drop proofErr wraps) ---
$old1 = "** Please have an outline in RMarkdown for your proejct with the topic, data source, and an overview for what you plan to do for your project as well as a timeline. Identify who (if anyone) you are working with and how you plan to divide the work. What question do you plan to answer?**"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2: merge "Devika Noir, Julia Smadja, and Dominic Thomas" ---
$old2 = "Devika Noir, Julia Smadja, and Dominic Thomas"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Change 3: merge "D.C. United has hired DJD to analyze their futbol match data so that they may learn which " ---
$old3 = "D.C. United has hired DJD to analyze their futbol match data so that they may learn which "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# --- Change 4: merge the "ting the outcome...this project." run (also removes the _GoBack bookmark that sat inside it) ---
$old4 = "ting the outcome of a game given certain conditions and test the models validity. Finally, DJD plans to develop an interactive and visualisation app so that D.C. United can continue to draw insights from the data at the conclusion of this project."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2) | Out-Null

# --- Change 5: merge "WUnderground: Weather History" ---
$old5 = "WUnderground: Weather History"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# --- Change 6: append " DONE" after the timeline's first bullet, and move the _GoBack bookmark there ---
$r = $d.Content
$r.Find.Execute("Data Ingestion/Cleaning - Feb 1 - Mar 1") | Out-Null
$r.Collapse(0)
$r.InsertAfter(" DONE")
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

Write-Output "ok"
